$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Common")
try {
  $ws.TopLeftCell = "A130"
  Write-Host "set ws.TopLeftCell ok"
} catch {
  Write-Host ("err: " + $_)
}
